{"js": "// Office.js (Word JavaScript API) edit script\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst REVISED_ABSTRACT = \"A medicina baseada em provas permite ao m\\u00e9dico avaliar a rela\\u00e7\\u00e3o risco-benef\\u00edcio de um tratamento atrav\\u00e9s da defini\\u00e7\\u00e3o e dos dados. As escolhas baseadas no risco podem ser feitas pelo m\\u00e9dico utilizando informa\\u00e7\\u00f5es diferentes. Com a emerg\\u00eancia de novas tecnologias, uma grande quantidade de dados \\u00e9 registada, oferecendo perspectivas interessantes com a aprendizagem de m\\u00e1quinas para a an\\u00e1lise preditiva de dados. A aprendizagem mec\\u00e2nica \\u00e9 um conjunto de m\\u00e9todos que processam dados para modelar um problema de aprendizagem. Algoritmos de aprendizagem supervisionada por m\\u00e1quinas consistem em utilizar dados anotados para construir o modelo. Esta categoria permite resolver problemas de an\\u00e1lise de dados de previs\\u00e3o. Neste artigo, detalhamos a utiliza\\u00e7\\u00e3o de algoritmos de aprendizagem supervisionada de m\\u00e1quinas para a previs\\u00e3o de problemas de an\\u00e1lise de dados em medicina. No campo m\\u00e9dico, os dados podem ser divididos em duas categorias: imagens m\\u00e9dicas e outros dados. Para uma maior brevidade, a nossa revis\\u00e3o trata de qualquer tipo de dados m\\u00e9dicos, excluindo imagens. Neste artigo, oferecemos uma discuss\\u00e3o em torno de quatro abordagens de aprendizagem supervisionada por m\\u00e1quinas: abordagens baseadas na informa\\u00e7\\u00e3o, baseadas na semelhan\\u00e7a, baseadas na probabilidade e baseadas em erros. Cada m\\u00e9todo \\u00e9 ilustrado com exemplos detalhados de medicina cardiovascular e nuclear. A nossa revis\\u00e3o mostra que o conjunto de modelos (ME) e a m\\u00e1quina vectorial de suporte (SVM) s\\u00e3o os m\\u00e9todos mais populares. SVM, ME e redes neurais artificiais conduzem frequentemente a melhores resultados do que os dados por outros algoritmos. Nos pr\\u00f3ximos anos, mais estudos, mais dados, mais ferramentas e mais m\\u00e9todos ser\\u00e3o, com certeza, propostos.\";\nconst CONCLUSAO_LABEL = \"Conclus\\u00e3o: \";\nconst CONCLUSION_TEXT = \"As aplica\\u00e7\\u00f5es de aprendizagem de m\\u00e1quinas est\\u00e3o a crescer no campo da medicina. Nos pr\\u00f3ximos anos, mais estudos, mais dados, mais ferramentas e mais m\\u00e9todos ser\\u00e3o, com certeza, propostos.\";\n\n// The document currently has 12 paragraphs (index 0-11). Paragraph 11 is the\n// long Portuguese abstract paragraph that gets reworded and followed by a new\n// \"Conclus\u00e3o\" section. Every paragraph in the document becomes justified.\n\nconst count = paragraphs.items.length;\n\n// 1) Justify every existing paragraph.\nfor (let i = 0; i < count; i++) {\n  paragraphs.items[i].alignment = Word.Alignment.justified;\n}\nawait context.sync();\n\n// 2) Replace the text of the final (abstract) paragraph with the revised wording.\nconst lastParagraph = paragraphs.items[count - 1];\nconst lastRange = lastParagraph.getRange();\nlastRange.insertText(REVISED_ABSTRACT, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Append the new \"Conclus\u00e3o\" block after the abstract paragraph:\n//    (empty) -> \"Conclus\u00e3o: \" -> (empty) -> conclusion text -> (empty)\nlet anchor = lastParagraph;\nanchor = anchor.insertParagraph(null, Word.InsertLocation.after);\nanchor.alignment = Word.Alignment.justified;\n\nanchor = anchor.insertParagraph(CONCLUSAO_LABEL, Word.InsertLocation.after);\nanchor.alignment = Word.Alignment.justified;\n\nanchor = anchor.insertParagraph(null, Word.InsertLocation.after);\nanchor.alignment = Word.Alignment.justified;\n\nanchor = anchor.insertParagraph(CONCLUSION_TEXT, Word.InsertLocation.after);\nanchor.alignment = Word.Alignment.justified;\n\nanchor = anchor.insertParagraph(null, Word.InsertLocation.after);\nanchor.alignment = Word.Alignment.justified;\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script (PowerShell-style)\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Justify every paragraph currently in the document.\nforeach ($p in $d.Paragraphs) {\n    $p.Format.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n}\n\n# 2) Replace the text of the last paragraph (the long Portuguese abstract)\n#    with the revised wording. Assigning to Range.Text keeps the paragraph's\n#    own formatting (pPr/rPr, language, etc.) intact.\n$count = $d.Paragraphs.Count\n$abstractParagraph = $d.Paragraphs.Item($count)\n$abstractParagraph.Range.Text = \"A medicina baseada em provas permite ao m\u00e9dico avaliar a rela\u00e7\u00e3o risco-benef\u00edcio de um tratamento atrav\u00e9s da defini\u00e7\u00e3o e dos dados. As escolhas baseadas no risco podem ser feitas pelo m\u00e9dico utilizando informa\u00e7\u00f5es diferentes. Com a emerg\u00eancia de novas tecnologias, uma grande quantidade de dados \u00e9 registada, oferecendo perspectivas interessantes com a aprendizagem de m\u00e1quinas para a an\u00e1lise preditiva de dados. A aprendizagem mec\u00e2nica \u00e9 um conjunto de m\u00e9todos que processam dados para modelar um problema de aprendizagem. Algoritmos de aprendizagem supervisionada por m\u00e1quinas consistem em utilizar dados anotados para construir o modelo. Esta categoria permite resolver problemas de an\u00e1lise de dados de previs\u00e3o. Neste artigo, detalhamos a utiliza\u00e7\u00e3o de algoritmos de aprendizagem supervisionada de m\u00e1quinas para a previs\u00e3o de problemas de an\u00e1lise de dados em medicina. No campo m\u00e9dico, os dados podem ser divididos em duas categorias: imagens m\u00e9dicas e outros dados. Para uma maior brevidade, a nossa revis\u00e3o trata de qualquer tipo de dados m\u00e9dicos, excluindo imagens. Neste artigo, oferecemos uma discuss\u00e3o em torno de quatro abordagens de aprendizagem supervisionada por m\u00e1quinas: abordagens baseadas na informa\u00e7\u00e3o, baseadas na semelhan\u00e7a, baseadas na probabilidade e baseadas em erros. Cada m\u00e9todo \u00e9 ilustrado com exemplos detalhados de medicina cardiovascular e nuclear. A nossa revis\u00e3o mostra que o conjunto de modelos (ME) e a m\u00e1quina vectorial de suporte (SVM) s\u00e3o os m\u00e9todos mais populares. SVM, ME e redes neurais artificiais conduzem frequentemente a melhores resultados do que os dados por outros algoritmos. Nos pr\u00f3ximos anos, mais estudos, mais dados, mais ferramentas e mais m\u00e9todos ser\u00e3o, com certeza, propostos.\"\n\n# 3) Append the new \"Conclus\u00e3o\" block right after the abstract paragraph:\n#    (empty) -> \"Conclus\u00e3o: \" -> (empty) -> conclusion text -> (empty)\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n$anchor.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n$anchor.Format.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n$anchor.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n$anchor.Range.Text = \"Conclus\u00e3o: \"\n$anchor.Format.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n$anchor.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n$anchor.Format.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n$anchor.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n$anchor.Range.Text = \"As aplica\u00e7\u00f5es de aprendizagem de m\u00e1quinas est\u00e3o a crescer no campo da medicina. Nos pr\u00f3ximos anos, mais estudos, mais dados, mais ferramentas e mais m\u00e9todos ser\u00e3o, com certeza, propostos.\"\n$anchor.Format.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n$anchor.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$anchor = $d.Paragraphs.Item($count)\n$anchor.Format.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify\n\nWrite-Output \"done\"\n"}
